# "1des fpoo aula05 correcao"
# Fills in the attendance ("chamada") column for the class held on the date
# in column V (2022-08-16) for every student row (3..51), mirroring the
# existing P/F (Presente/Falta) marks already present in column U, with a
# handful of manual corrections, then sets up the same freeze-panes /
# selection view state that was captured in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number (1-based, matching the worksheet) -> attendance mark for
# column V ("P" = Presente, "F" = Falta).
$attendance = [ordered]@{
  3  = "P"
  4  = "P"
  5  = "P"
  6  = "P"
  7  = "P"
  8  = "F"
  9  = "P"
  10 = "P"
  11 = "P"
  12 = "P"
  13 = "P"
  14 = "P"
  15 = "P"
  16 = "P"
  17 = "P"
  18 = "F"
  19 = "F"
  20 = "P"
  21 = "P"
  22 = "F"
  23 = "P"
  24 = "P"
  25 = "P"
  26 = "P"
  27 = "F"
  28 = "P"
  29 = "F"
  30 = "F"
  31 = "P"
  32 = "P"
  33 = "P"
  34 = "F"
  35 = "P"
  36 = "P"
  37 = "P"
  38 = "P"
  39 = "F"
  40 = "P"
  41 = "F"
  42 = "P"
  43 = "F"
  44 = "F"
  45 = "P"
  46 = "F"
  47 = "F"
  48 = "P"
  49 = "P"
  50 = "P"
  51 = "P"
}

foreach ($row in $attendance.Keys) {
    $ws.Cells.Item($row, 22).Value = $attendance[$row]
}

# Recreate the saved view: freeze the first two columns (A:B) and the first
# two rows (1:2), then leave the active selection on V8 (the class-05 column
# for the row the teacher was editing), matching the panes/selection saved
# in the sheet view.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("C3").Select()
$win.FreezePanes = $true
$ws.Range("V8").Select()
